$wb = $excel.ActiveWorkbook

# --- preparation_medium: add 4 new rows (Alpha-MEM, Modified Davidson's Fixative,
#     Growth media, Lysis buffer) ---
$wsPrepMedium = $wb.Worksheets.Item("preparation_medium")

$wsPrepMedium.Rows.Item(11).Insert()
$wsPrepMedium.Range("A11").Value = "Alpha-MEM"
$wsPrepMedium.Range("B11").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000371"

$wsPrepMedium.Rows.Item(18).Insert()
$wsPrepMedium.Range("A18").Value = "Modified Davidson's Fixative"
$wsPrepMedium.Range("B18").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000370"

$wsPrepMedium.Rows.Item(24).Insert()
$wsPrepMedium.Range("A24").Value = "Growth media"
$wsPrepMedium.Range("B24").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000372"

$wsPrepMedium.Range("A30").Value = "Lysis buffer"
$wsPrepMedium.Range("B30").Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C178573"

# --- storage_medium: add 1 new row (Formic acid in water) ---
$wsStorageMedium = $wb.Worksheets.Item("storage_medium")

$wsStorageMedium.Rows.Item(6).Insert()
$wsStorageMedium.Range("A6").Value = "Formic acid in water"
$wsStorageMedium.Range("B6").Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C83719"

# --- Sample Suspension: widen the data validation dropdown ranges to cover
#     the newly added rows in preparation_medium and storage_medium ---
$wsMain = $wb.Worksheets.Item("Sample Suspension")

$wsMain.Range("I2:I1001").Validation.Formula1 = "'preparation_medium'!`$A`$1:`$A`$30"
$wsMain.Range("M2:M1001").Validation.Formula1 = "'storage_medium'!`$A`$1:`$A`$21"

# --- .metadata: bump the pav:createdOn timestamp ---
$wsMetadata = $wb.Worksheets.Item(".metadata")
$wsMetadata.Range("C2").Value = "2024-04-20T17:32:34-07:00"
